$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Row 1 (table row index 1)
Replace-Text "54×48=2592" "94×29=2726"
Replace-Text "54×45=2430" "83×96=7968"
Replace-Text "36×96=3456" "96×88=8448"
Replace-Text "22×22=484" "28×33=924"
Replace-Text "90×19=1710" "81×85=6885"

# Row 5
Replace-Text "31×34=1054" "57×30=1710"
Replace-Text "48×82=3936" "36×51=1836"
Replace-Text "97×11=1067" "71×51=3621"
Replace-Text "75×59=4425" "17×44=748"
Replace-Text "50×53=2650" "92×87=8004"

# Row 10: cell contents are reshuffled (two old cells drop out, two new ones
# appear, and the surviving "40x44=1760" value moves from column 3 to
# column 5). Net effect per-column, addressed directly via the table model
# to avoid any ambiguity from the repeated "40x44=1760" text.
$tbl = $d.Tables.Item(1)
$row9 = $tbl.Rows.Item(10)
$row9.Cells.Item(1).Range.Text = "61×86=5246"
$row9.Cells.Item(2).Range.Text = "24×23=552"
$row9.Cells.Item(3).Range.Text = "43×72=3096"
$row9.Cells.Item(4).Range.Text = "19×23=437"
$row9.Cells.Item(5).Range.Text = "40×44=1760"

# Row 15
Replace-Text "61×37=2257" "88×59=5192"
Replace-Text "33×85=2805" "81×49=3969"
Replace-Text "20×18=360" "39×65=2535"
Replace-Text "17×67=1139" "12×15=180"
Replace-Text "59×82=4838" "36×69=2484"

# Row 20
Replace-Text "82×39=3198" "66×46=3036"
Replace-Text "20×96=1920" "44×18=792"
Replace-Text "30×43=1290" "47×38=1786"
Replace-Text "52×62=3224" "52×11=572"
Replace-Text "92×13=1196" "63×67=4221"
